# JS-Frameworks-Self-Evaluation-Protocol.xlsx
# "project "finished" self evaluation updated"
#
# - Fill in the "Score" (C) column with the grader's actual scores.
# - Clear out the old "Comments" (E) column values (no longer used as a
#   data column once the real scores live in C).
# - C8/D8/E8 previously carried a leftover date number-format (from when
#   this row used to hold a date) - reset them to General before writing
#   the new "Days commit" / "Number of commits" values.
# - Total Score (C44) becomes a real SUM formula over the new score column.
# - Misc cosmetic: refresh window selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8/9: GitHub activity numbers -------------------------------------
# These two cells used to be styled/number-formatted as dates (numFmtId=16,
# "d-mmm") - reset to General so the new numeric/"infinity" values display
# correctly instead of turning into serial-date text.
$ws.Range("C8").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "General"

$infinity = [char]0x221E

$ws.Range("C8").Value = 6
$ws.Range("D8").Value = $infinity
$ws.Range("E8").Value = $null

$ws.Range("C9").Value = 52
$ws.Range("D9").Value = $infinity
$ws.Range("E9").Value = $null

# --- Basic Options score table (rows 11-34) --------------------------------
# Score (C) / Comments (E) pairs: fill in C with the real score, clear E.
$scores = @{
    11 = 5
    12 = 20
    13 = 5
    14 = 10
    15 = 5
    16 = 10
    17 = 2
    18 = 6
    19 = 10
    20 = 5
    21 = 9
    22 = 9
    23 = 5
    24 = 5
    25 = 9
    26 = 9
    27 = 9
    30 = 10
    31 = 5
    32 = 5
    33 = 5
    34 = 8
}
foreach ($r in $scores.Keys) {
    $ws.Range("C$r").Value = $scores[$r]
    $ws.Range("E$r").Value = $null
}

# Rows 28/29 keep an empty score, but still lose their old Comments value.
$ws.Range("E28").Value = $null
$ws.Range("E29").Value = $null

# --- Advanced Options score table (rows 36-43) ------------------------------
$advScores = @{
    36 = 10
    37 = 5
    38 = 10
    39 = 5
    42 = 20
    43 = 10
}
foreach ($r in $advScores.Keys) {
    $ws.Range("C$r").Value = $advScores[$r]
    $ws.Range("E$r").Value = $null
}

# --- Total score becomes a real formula ------------------------------------
$ws.Range("C44").Formula = "=SUM(C6:C43)"

# --- Cosmetic: selection / scroll position ----------------------------------
$ws.Range("E41").Select()
